# Update Automation TestCase & Update Excel Data
#
# The "Users" worksheet has a table (Table5, A1:B2) listing login
# credentials used by the Selenium test-suite. A new test account
# (hthuy / 123) was added as a third table row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Grow the Users table by one row (expands ref A1:B2 -> A1:B3 and
# updates the table/autofilter ranges automatically).
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Copy down the formatting from the row above (row 2) onto the new
# row 3, same as Excel does when you extend a table by typing into
# the row right below it.
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122)

# Fill in the new credentials.
$ws.Cells.Item(3, 1).Value = "hthuy"
$ws.Cells.Item(3, 2).Value = "123"

# Leave the selection where the author last clicked while editing.
$ws.Activate()
$ws.Range("D11").Select()
